$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear header cell A1 and remove the bold/bordered/centered style from row 1 ---
$ws.Range("A1").Value = ""
$ws.Range("A1:O1").Style = "Normal"

# --- Update numeric data in rows 3-7 (corrected pre/post/total fixation data) ---

# Row 3 - Revisit count
$ws.Range("B3").Value = 14
$ws.Range("C3").Value = 66
$ws.Range("D3").Value = 21
$ws.Range("E3").Value = 25
$ws.Range("F3").Value = 21
$ws.Range("G3").Value = 2
$ws.Range("I3").Value = 8
$ws.Range("J3").Value = 36
$ws.Range("K3").Value = 18
$ws.Range("L3").Value = 77
$ws.Range("M3").Value = 3
$ws.Range("N3").Value = 1

# Row 4 - Fixation count
$ws.Range("B4").Value = 26
$ws.Range("C4").Value = 257
$ws.Range("D4").Value = 41
$ws.Range("E4").Value = 44
$ws.Range("F4").Value = 41
$ws.Range("G4").Value = 5
$ws.Range("I4").Value = 12
$ws.Range("J4").Value = 63
$ws.Range("K4").Value = 31
$ws.Range("L4").Value = 365
$ws.Range("M4").Value = 4
$ws.Range("N4").Value = 2

# Row 5 - Dwell time (ms)
$ws.Range("B5").Value = 13580.45
$ws.Range("C5").Value = 78274.68
$ws.Range("D5").Value = 19019.04
$ws.Range("E5").Value = 16140.64
$ws.Range("F5").Value = 19019.04
$ws.Range("G5").Value = 3686.74
$ws.Range("I5").Value = 7824.8
$ws.Range("J5").Value = 19803.06
$ws.Range("K5").Value = 9585.95
$ws.Range("L5").Value = 111690.9
$ws.Range("M5").Value = 1468.02
$ws.Range("N5").Value = 1117.65

# Row 6 - Dwell time (%)
$ws.Range("B6").Value = 4.43
$ws.Range("C6").Value = 25.53
$ws.Range("D6").Value = 6.2
$ws.Range("E6").Value = 5.27
$ws.Range("F6").Value = 6.2
$ws.Range("G6").Value = 1.2
$ws.Range("I6").Value = 2.55
$ws.Range("J6").Value = 6.46
$ws.Range("K6").Value = 3.13
$ws.Range("L6").Value = 36.43
$ws.Range("M6").Value = 0.48
$ws.Range("N6").Value = 0.36
$ws.Range("O6").Value = 0.05

# Row 7 - Fixation duration (ms)
$ws.Range("B7").Value = 522.33
$ws.Range("C7").Value = 304.57
$ws.Range("D7").Value = 463.88
$ws.Range("E7").Value = 366.83
$ws.Range("F7").Value = 463.88
$ws.Range("G7").Value = 737.35
$ws.Range("I7").Value = 652.07
$ws.Range("J7").Value = 314.33
$ws.Range("K7").Value = 309.22
$ws.Range("L7").Value = 306
$ws.Range("M7").Value = 367
$ws.Range("N7").Value = 558.82

# --- Remove now-empty trailing rows 10-14 ---
$ws.Range("A10:O14").EntireRow.Delete()
